# close #124: Add function to validate sequential codes in spreadsheet
#
# The header row (row 1) codes, C1:S1, used to be based on the old
# 5000-series ids (e.g. "5000-2015", "5003-2030-O", ...). They are
# renamed here to plain sequential codes (1, 2, 3, ... 9) so a
# "validate sequential codes" check can be run against them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old header -> new header, column by column (row 1, columns C..S)
$ws.Range("C1").Value = "1-2015"
$ws.Range("D1").Value = "2-2015"
$ws.Range("E1").Value = "2-2030-O"
$ws.Range("F1").Value = "2-2050-O"
$ws.Range("G1").Value = "2-2030-P"
$ws.Range("H1").Value = "2-2050-P"
$ws.Range("I1").Value = "3-2015"
$ws.Range("J1").Value = "4-2015"
$ws.Range("K1").Value = "5-2015"
$ws.Range("L1").Value = "5-2030-O"
$ws.Range("M1").Value = "5-2050-O"
$ws.Range("N1").Value = "5-2030-P"
$ws.Range("O1").Value = "5-2050-P"
$ws.Range("P1").Value = "6-2015"
$ws.Range("Q1").Value = "7-2015"
$ws.Range("R1").Value = "8-2015"
$ws.Range("S1").Value = "9-2015"

# Small numeric corrections found while validating the sequential codes.
$ws.Range("J11").Value = 0.345198133949
$ws.Range("G14").Value = 0.867444431328
$ws.Range("F15").Value = 0.274029954469

# Move the active selection.
$ws.Range("B22").Select()
